$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F2:F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1042
$ws1.Range("F3").Value = 234
$ws1.Range("F4").Value = 2559
$ws1.Range("F5").Value = 42
$ws1.Range("F6").Value = 552

# Sheet "全部类型" (All Types) - update F4:F8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1042
$ws4.Range("F5").Value = 234
$ws4.Range("F6").Value = 2559
$ws4.Range("F7").Value = 42
$ws4.Range("F8").Value = 552
